$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "remove-uneeded-dependencies" (row 32) is no longer shipped as a separate
# script entry - delete its row entirely. This shifts soundconverter/steam/etc.
# up by one row (A1:D49 -> A1:D48) and drops the now-unused shared string.
$null = $ws.Rows(32).Delete()

# soundconverter (now row 32) and steam (now row 33) are supported on Debian -
# mark their "Debian" column (D) as Supported, matching the green formatting
# already used elsewhere in the column.
$null = $ws.Range("C32").Copy()
$null = $ws.Range("D32").PasteSpecial(-4122)
$null = $ws.Range("C33").Copy()
$null = $ws.Range("D33").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

$null = $ws.Range("D33").Select()
